$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Target change (per the commit "The Dragon Boat Festival"):
#   - Para 1 "2022年6月1日星期三" keeps its text, but loses the _GoBack bookmark.
#   - A brand-new paragraph is inserted right after it, containing the text
#     that used to live in (old) paragraph 2: "多云，今天是六一儿童节，应收款项的一天".
#   - The old paragraph 2 keeps its own paragraph (now the 3rd one) but its
#     text becomes "哈无DAU无敌摩擦的温暖的", and the _GoBack bookmark moves to
#     the end of that paragraph (after the run, before the paragraph mark).
# ------------------------------------------------------------------

# Step 1: Update the original second paragraph's text to the new wording.
# Do this BEFORE inserting the duplicated paragraph below, so the Find only
# matches the one (still unique) occurrence of the old text.
$d.Content.Find.Execute("多云，今天是六一儿童节，应收款项的一天", $true, $false, $false, $false, $false, $true, 1, $false, "哈无DAU无敌摩擦的温暖的", 2)

# Step 2: Insert a new paragraph right after paragraph 1 and give it the text
# that originally lived in paragraph 2 (now duplicated here).
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "多云，今天是六一儿童节，应收款项的一天"

# Step 3: Move the _GoBack bookmark from the end of paragraph 1 to the end of
# paragraph 3 (the renumbered old paragraph 2), landing after its run and
# before the paragraph mark.
#
# Note: adding a bookmark exactly at a position that sits right before a
# paragraph mark (paragraph.Range.End - 1) confuses this host's bookmark
# placement, so a one-character scratch placeholder is used to borrow a safe
# spot, then removed once the bookmark is anchored correctly.
$p3 = $d.Paragraphs(3)
$scratchPos = $p3.Range.End - 1
$scratchRange = $d.Range($scratchPos, $scratchPos)
$scratchRange.InsertBefore("Z")

$p3b = $d.Paragraphs(3)
$safePos = $p3b.Range.End - 2
$bmRange = $d.Range($safePos, $safePos)

$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $bmRange)

$bm = $d.Bookmarks("_GoBack")
$scratchChar = $d.Range($bm.End, $bm.End + 1)
$scratchChar.Delete()

Write-Output "Done."
